$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: NES, pvalue, p.adjust, qvalues
$ws.Range("E2").Value = -2.66841154794926493921
$ws.Range("F2").Value = 0.00006919818492609946
$ws.Range("G2").Value = 0.00166075643822638681
$ws.Range("H2").Value = 0.00160248428249914511

# Row 3: NES, pvalue, p.adjust, qvalues
$ws.Range("E3").Value = -2.22749824646499927638
$ws.Range("F3").Value = 0.0019636350033408142
$ws.Range("G3").Value = 0.02356362004008977218
$ws.Range("H3").Value = 0.02273682635447258835
